$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("singleInputData")
$ws2 = $wb.Worksheets.Item("doubleInputData")

# Add the new test-case row to each sheet. The order in which NEW string
# values are first written controls their position in the shared-strings
# table, so this mirrors the order needed to reproduce the target file:
# Faile(8), Fail(9), Failed for Raport Purpoose(10), Raport(11)
$ws1.Range("A6").Value = "Faile"
$ws2.Range("A6").Value = "Fail"
$ws1.Range("B6").Value = "Failed for Raport Purpoose"
$ws2.Range("B6").Value = "Raport"
$ws2.Range("C6").Value = "Failed for Raport Purpoose"

# Widen the columns that now hold the longer strings.
$ws1.Columns.Item(2).ColumnWidth = 23
$ws2.Columns.Item(3).ColumnWidth = 24.17

# Update the remembered selection on the first sheet, then move to/select
# on the second sheet last, so it becomes the active (tabSelected) sheet.
$ws1.Activate()
$ws1.Range("B6").Select()

$ws2.Activate()
$ws2.Range("B6").Select()
